$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77. This shifts the existing rows 77-161
# down to 78-162 (format of the row above, including the date style on
# column D, carries down onto the freshly inserted row automatically).
$ws.Rows(77).Insert()

# Populate the newly-inserted row 77 with the new data record. The
# non-numeric/non-date columns mirror what ends up in (old) row 77 -
# now shifted to row 78 - since the diff only shows D/M/N/O/P/S changing
# for this record.
$ws.Range("A77").Value = 9
$ws.Range("B77").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C77").Value = "Metropolitana"
$ws.Range("D77").Value = 45049
$ws.Range("E77").Value = 13
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100101
$ws.Range("H77").Value = "Berries"
$ws.Range("I77").Value = 100101004
$ws.Range("J77").Value = "Frambuesa"
$ws.Range("K77").Value = "Sin especificar"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 330
$ws.Range("N77").Value = 7000
$ws.Range("O77").Value = 7500
$ws.Range("P77").Value = 7227
$ws.Range("Q77").Value = "`$/bandeja 2 kilos"
$ws.Range("R77").Value = "Provincia de Linares"
$ws.Range("S77").Value = 3614
$ws.Range("T77").Value = 2
